# "Generate Report for Handoff" - regenerate the localization status report
# with a freshly generated source markdown file (new GUID) and xliff
# handoff files, clearing out the stale handback info for the new round.

$wb = $excel.ActiveWorkbook

$oldGuid = "f5f8611f-a9cd-411b-9779-ac65a62fc8c0"
$newGuid = "c0380e28-1ff3-4947-9f21-58095cb79a70"

$newFileName      = "$newGuid.md"
$newPathAndName   = "e2e\$newGuid.md"
$newZhXlf         = "$newGuid.391c6050751e83464e0aa9a79235d57def9efd4d.zh-cn.xlf"
$newDeXlf         = "$newGuid.391c6050751e83464e0aa9a79235d57def9efd4d.de-de.xlf"

$hoGenerateDate   = "2016-09-07 17:28:25"
$zhHandoffDate    = "2016-09-07 17:28:19"
$deHandoffDate    = $hoGenerateDate
$emptyHandback    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $hoGenerateDate

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = $newPathAndName
    }
}

$wsOverview.Columns.Item(1).ColumnWidth = 39.1667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("K2").Value = $emptyHandback

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFileName
    }
    if ($addr -eq '$I$2') {
        $h.Delete()
    }
}

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("J2").Style = "Normal"

$wsZh.Columns.Item(1).ColumnWidth = 39.1667
$wsZh.Columns.Item(9).ColumnWidth = 17.8333
$wsZh.Columns.Item(10).ColumnWidth = 20.8333

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $deHandoffDate
$wsDe.Range("K2").Value = $emptyHandback

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFileName
    }
    if ($addr -eq '$I$2') {
        $h.Delete()
    }
}

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("J2").Style = "Normal"

$wsDe.Columns.Item(1).ColumnWidth = 39.1667
$wsDe.Columns.Item(9).ColumnWidth = 17.8333
$wsDe.Columns.Item(10).ColumnWidth = 20.8333
